$wb = $excel.ActiveWorkbook

# Column F holds the "want to go" attendance counts. This data refresh
# (gh-pages regeneration) bumps those counts on both the "展览" sheet and
# the combined "全部类型" sheet, which duplicates the same rows.
$targets = @(
    @{ Row = 2;  Value = 1044 },
    @{ Row = 6;  Value = 132 },
    @{ Row = 10; Value = 5110 },
    @{ Row = 11; Value = 4740 }
)

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($t in $targets) {
        $ws.Cells.Item($t.Row, 6).Value = $t.Value
    }
}
